$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: B6, C6, E6 were stored as text ("56348.0", "-318.0", "-363.0") but
# should become real numbers (56348, -318, -363). D6/F6/G6 stay as-is.
$ws.Range("B6").Value = 56348
$ws.Range("C6").Value = -318
$ws.Range("E6").Value = -363

# Row 7 (new): values look numeric/date/percent-like, so a plain .Value
# assignment would make Excel auto-convert them into a real date/number.
# Prefix with an apostrophe to force literal text entry (same as typing
# '2022-01-03 into the cell), then reset the style back to Normal so we
# don't leave the quote-prefix style behind.
$ws.Range("A7").Value = "'2022-01-03"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "'56348.0"
$ws.Range("B7").Style = "Normal"

$ws.Range("C7").Value = "'-1065.0"
$ws.Range("C7").Style = "Normal"

$ws.Range("D7").Value = "'-1.89%"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "'0"
$ws.Range("E7").Style = "Normal"

# F7 / G7 stay blank/empty text, matching the blank trailing columns of
# row 6. A lone apostrophe is an empty text entry (quote-prefix with no
# characters after it), giving an empty-string text cell instead of a
# truly-blank one; reset the style afterwards so no quote-prefix style
# carries over.
$ws.Range("F7").Value = "'"
$ws.Range("F7").Style = "Normal"

$ws.Range("G7").Value = "'"
$ws.Range("G7").Style = "Normal"
